$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '22.397.72'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +0.31%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.569.94'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +0.07%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.24%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.003'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.26%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '289.72'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.01%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3744'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.25%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '49.30'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.53%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3360'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -1.57%  '

$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07448'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -2.94%  '

$ws.Range("B11").Value = 'Polygon'
$ws.Range("C11").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.127'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -3.48%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.003'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +0.25%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.92'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -2.17%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.900'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -2.16%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.861'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -1.28%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.568.61'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -0.17%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001115'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -2.11%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '89.02'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -1.37%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06684'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.55%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.003'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +0.26%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.151'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -1.65%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '16.15'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -2.93%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.86'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.95%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '22.397.33'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.32%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.367'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -1.04%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.517'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -9.80%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.97'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.89%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '147.22'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +1.27%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.997'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +0.16%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '124.62'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -0.90%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.742.19'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -0.02%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9995'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -2.38%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.964'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -2.35%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.883'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -5.65%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.695'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -4.37%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.08413'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -1.28%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.369'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +3.35%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02441'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -3.94%  '

$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06451'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +0.63%  '

$ws.Range("B40").Value = 'Algorand'
$ws.Range("C40").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2242'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -3.71%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.372'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -3.40%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.32'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -3.93%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.6197'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -3.52%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.003'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +0.31%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.93'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -2.15%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.808'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +1.31%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5785'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -3.57%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.053'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -2.20%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '125.27'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +0.48%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.227'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -3.67%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07297'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +0.23%  '
